$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date for each row of data
# (rows 2-16). The automatic update bumps this date from 2023-10-13
# (serial 45212) to 2023-10-22 (serial 45221) for every row, leaving the
# existing date formatting/style untouched.
for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
